$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 30539.766
$ws.Range("I74").Value = 42998.086
$ws.Range("J74").Value = 4490.5454
$ws.Range("K74").Value = 42998.086
$ws.Range("L74").Value = 4490.5454
$ws.Range("M74").Value = -42062.086
$ws.Range("N74").Value = -6362.5454
$ws.Range("H77").Value = 30539.766
$ws.Range("I77").Value = 42998.086
$ws.Range("J77").Value = 4490.5454
$ws.Range("K77").Value = 214990.43
$ws.Range("L77").Value = 22452.727
$ws.Range("M77").Value = -210310.43
$ws.Range("N77").Value = -31812.727
$ws.Range("H112").Value = 2505
$ws.Range("J112").Value = 2594.4119
$ws.Range("L112").Value = 7783.2357
$ws.Range("N112").Value = -9999.235700000001
$ws.Range("H113").Value = 3725
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").Value = $null
$ws.Range("H129").Value = 1092.2703
$ws.Range("I129").Value = 462.2
$ws.Range("K129").Value = 1386.6
$ws.Range("M129").Value = 3613.4
$ws.Range("H137").Value = 1817.5714
$ws.Range("I137").Value = 1112.5518
$ws.Range("J137").Value = 2574.8147
$ws.Range("K137").Value = 3337.6554
$ws.Range("L137").Value = 7724.4441
$ws.Range("M137").Value = -787.6553999999996
$ws.Range("N137").Value = -12824.4441
$ws.Range("H138").Value = 3941.7568
$ws.Range("I138").Value = 2537.3684
$ws.Range("J138").Value = 4426.909
$ws.Range("K138").Value = 7612.1052
$ws.Range("L138").Value = 13280.727
$ws.Range("M138").Value = -2472.1052
$ws.Range("N138").Value = -23560.727

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12497.125
$ws.Range("I32").Value = 12063.6
$ws.Range("K32").Value = 12063.6
$ws.Range("M32").Value = -11776.6
$ws.Range("H33").Value = 24000
$ws.Range("I33").Value = 24000
$ws.Range("K33").Value = 24000
$ws.Range("M33").Value = -23671
$ws.Range("H132").Value = 4879.6904
$ws.Range("I132").Value = 4962.0605
$ws.Range("J132").Value = 4577.6665
$ws.Range("K132").Value = 14886.1815
$ws.Range("L132").Value = 13732.9995
$ws.Range("M132").Value = -12356.1815
$ws.Range("N132").Value = -18792.9995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 7290
$ws.Range("I22").Value = 10086
$ws.Range("J22").Value = 300
$ws.Range("K22").Value = 10086
$ws.Range("L22").Value = 300
$ws.Range("M22").Value = -9913
$ws.Range("N22").Value = -646
$ws.Range("H105").Value = 9527069
$ws.Range("I105").Value = 12990047
$ws.Range("J105").Value = 3880.25
$ws.Range("K105").Value = 12990047
$ws.Range("L105").Value = 3880.25
$ws.Range("M105").Value = -12988300
$ws.Range("N105").Value = -7374.25
$ws.Range("H107").Value = 28062.191
$ws.Range("I107").Value = 37487.133
$ws.Range("J107").Value = 4499.8335
$ws.Range("K107").Value = 37487.133
$ws.Range("L107").Value = 4499.8335
$ws.Range("M107").Value = -35567.133
$ws.Range("N107").Value = -8339.833500000001
$ws.Range("H134").Value = 3261.0527
$ws.Range("I134").Value = 2996.3635
$ws.Range("K134").Value = 8989.0905
$ws.Range("M134").Value = -6454.0905

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 5500
$ws.Range("H16").Value = 1435.8334
$ws.Range("I16").Value = 1207.3334
$ws.Range("J16").Value = 1664.3334
$ws.Range("K16").Value = 1207.3334
$ws.Range("L16").Value = 1664.3334
$ws.Range("M16").Value = -920.3334
$ws.Range("N16").Value = -2238.3334
$ws.Range("H58").Value = 1685388.1
$ws.Range("I58").Value = 2471042.5
$ws.Range("J58").Value = 1842.8572
$ws.Range("K58").Value = 2471042.5
$ws.Range("L58").Value = 1842.8572
$ws.Range("M58").Value = -2470839.5
$ws.Range("N58").Value = -2248.8572
$ws.Range("H70").Value = 48288
$ws.Range("J70").Value = 48288
$ws.Range("L70").Value = 48288
$ws.Range("N70").Value = -48918
$ws.Range("H73").Value = 48288
$ws.Range("J73").Value = 48288
$ws.Range("L73").Value = 48288
$ws.Range("N73").Value = -50472
$ws.Range("H104").Value = 23980
$ws.Range("J104").Value = 23980
$ws.Range("L104").Value = 23980
$ws.Range("N104").Value = -29222
$ws.Range("H107").Value = 1062.2
$ws.Range("I107").Value = 1062.2
$ws.Range("K107").Value = 1062.2
$ws.Range("M107").Value = 857.8
$ws.Range("H113").Value = 1435.8334
$ws.Range("I113").Value = 1207.3334
$ws.Range("J113").Value = 1664.3334
$ws.Range("K113").Value = 1207.3334
$ws.Range("L113").Value = 1664.3334
$ws.Range("M113").Value = 962.6666
$ws.Range("N113").Value = -6004.3334
$ws.Range("H136").Value = 1685388.1
$ws.Range("I136").Value = 2471042.5
$ws.Range("J136").Value = 1842.8572
$ws.Range("K136").Value = 7413127.5
$ws.Range("L136").Value = 5528.571599999999
$ws.Range("M136").Value = -7410577.5
$ws.Range("N136").Value = -10628.5716

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 4224.45
$ws.Range("I64").Value = 2990
$ws.Range("J64").Value = 5458.9
$ws.Range("K64").Value = 8970
$ws.Range("L64").Value = 16376.7
$ws.Range("M64").Value = -8700
$ws.Range("N64").Value = -16916.7
$ws.Range("H67").Value = 4224.45
$ws.Range("I67").Value = 2990
$ws.Range("J67").Value = 5458.9
$ws.Range("K67").Value = 8970
$ws.Range("L67").Value = 16376.7
$ws.Range("M67").Value = -8034
$ws.Range("N67").Value = -18248.7
$ws.Range("H68").Value = 1205.4166
$ws.Range("I68").Value = 881.46
$ws.Range("J68").Value = 1681.8235
$ws.Range("K68").Value = 2644.38
$ws.Range("L68").Value = 5045.470499999999
$ws.Range("M68").Value = -1833.38
$ws.Range("N68").Value = -6667.470499999999
$ws.Range("H71").Value = 1205.4166
$ws.Range("I71").Value = 881.46
$ws.Range("J71").Value = 1681.8235
$ws.Range("K71").Value = 7933.14
$ws.Range("L71").Value = 15136.4115
$ws.Range("M71").Value = -3877.14
$ws.Range("N71").Value = -23248.4115
$ws.Range("H113").Value = 590.37036
$ws.Range("I113").Value = 590.37036
$ws.Range("K113").Value = 1771.11108
$ws.Range("M113").Value = 398.8889199999999
$ws.Range("H132").Value = 1355.8572
$ws.Range("J132").Value = 1416.591
$ws.Range("L132").Value = 12749.319
$ws.Range("N132").Value = -17809.319

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2702.9583
$ws.Range("I80").Value = 2492.5
$ws.Range("J80").Value = 3334.3333
$ws.Range("K80").Value = 2492.5
$ws.Range("L80").Value = 3334.3333
$ws.Range("M80").Value = -1494.5
$ws.Range("N80").Value = -5330.3333
$ws.Range("H83").Value = 2702.9583
$ws.Range("I83").Value = 2492.5
$ws.Range("J83").Value = 3334.3333
$ws.Range("K83").Value = 12462.5
$ws.Range("L83").Value = 16671.6665
$ws.Range("M83").Value = -7470.5
$ws.Range("N83").Value = -26655.6665
$ws.Range("H102").Value = 3500.875
$ws.Range("I102").Value = 3487.8696
$ws.Range("J102").Value = 3800
$ws.Range("K102").Value = 3487.8696
$ws.Range("L102").Value = 3800
$ws.Range("M102").Value = -1865.8696
$ws.Range("N102").Value = -7044
$ws.Range("H113").Value = 1978.6522
$ws.Range("I113").Value = 1212.5
$ws.Range("J113").Value = 2387.2666
$ws.Range("K113").Value = 1212.5
$ws.Range("L113").Value = 2387.2666
$ws.Range("M113").Value = 957.5
$ws.Range("N113").Value = -6727.2666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 5000
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").Value = $null
$ws.Range("H27").Value = 5000
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").Value = $null
$ws.Range("H61").Value = 27436.75
$ws.Range("I61").Value = 51948.5
$ws.Range("K61").Value = 51948.5
$ws.Range("M61").Value = -51746.5
$ws.Range("H113").Value = 27436.75
$ws.Range("I113").Value = 51948.5
$ws.Range("K113").Value = 51948.5
$ws.Range("M113").Value = -49778.5
$ws.Range("H122").Value = 22503906
$ws.Range("I122").Value = 25002510
$ws.Range("J122").Value = 20005300
$ws.Range("K122").Value = 75007530
$ws.Range("L122").Value = 60015900
$ws.Range("M122").Value = -75005080
$ws.Range("N122").Value = -60020800
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").Value = $null
$ws.Range("H132").Value = 6220.273
$ws.Range("I132").Value = 6220.273
$ws.Range("K132").Value = 18660.819
$ws.Range("M132").Value = -16130.819

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 68457.47
$ws.Range("I81").Value = 81983.92999999999
$ws.Range("J81").Value = 5334
$ws.Range("K81").Value = 163967.86
$ws.Range("L81").Value = 10668
$ws.Range("M81").Value = -162906.86
$ws.Range("N81").Value = -12790
$ws.Range("H84").Value = 68457.47
$ws.Range("I84").Value = 81983.92999999999
$ws.Range("J84").Value = 5334
$ws.Range("K84").Value = 819839.2999999999
$ws.Range("L84").Value = 53340
$ws.Range("M84").Value = -814535.2999999999
$ws.Range("N84").Value = -63948
$ws.Range("H107").Value = 471.08334
$ws.Range("J107").Value = 775.5
$ws.Range("L107").Value = 2326.5
$ws.Range("N107").Value = -6166.5
$ws.Range("H113").Value = 452.9091
$ws.Range("I113").Value = 383.1
$ws.Range("K113").Value = 1149.3
$ws.Range("M113").Value = 1020.7
$ws.Range("H126").Value = 1497.5294
$ws.Range("I126").Value = 1397.2858
$ws.Range("J126").Value = 1965.3334
$ws.Range("K126").Value = 4191.857400000001
$ws.Range("L126").Value = 5896.0002
$ws.Range("M126").Value = -1721.857400000001
$ws.Range("N126").Value = -10836.0002
$ws.Range("H141").Value = 74118.89
$ws.Range("J141").Value = 74118.89
$ws.Range("L141").Value = 74118.89
$ws.Range("N141").Value = -84478.89
